$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "MME235-2016-17-(2)"
$ws.Range("C2").Value = "MME235-2014-15-(6)"
$ws.Range("D2").Value = "EEE267-2014-15-(1)"
$ws.Range("B3").Value = "ME221-2018-19-(1)"
$ws.Range("C3").Value = "ME221-2014-15-(4)"
$ws.Range("D3").Value = "EEE267-2017-18-(4)"
$ws.Range("B4").Value = "MME213-2015-16-(8)"
$ws.Range("D4").Value = "ME221-2015-16-(8)"
$ws.Range("B5").Value = "EEE267-2016-17-(3)"
$ws.Range("C5").Value = "EEE267-2015-16-(8)"
$ws.Range("D5").Value = "EEE267-2017-18-(8)"
$ws.Range("B6").Value = "MME213-2014-15-(12)"
$ws.Range("C6").Value = "MME213-2018-19-(4)"
$ws.Range("D6").Value = "ME221-2015-16-(5)"
$ws.Range("B7").Value = "EEE267-2018-19-(2)"
$ws.Range("C7").Value = "EEE267-2018-19-(8)"
$ws.Range("D7").Value = "MME213-2016-17-(6)"
$ws.Range("B8").Value = "ME221-2016-17-(3)"
$ws.Range("C8").Value = "MME213-2015-16-(9)"
$ws.Range("D8").Value = "MME213-2016-17-(9)"
$ws.Range("B9").Value = "ME221-2014-15-(7)"
$ws.Range("C9").Value = "EEE267-2015-16-(5)"
$ws.Range("D9").Value = "MME235-2017-18-(3)"
$ws.Range("B10").Value = "MME213-2017-18-(4)"
$ws.Range("C10").Value = "EEE267-2017-18-(1)"
$ws.Range("D10").Value = "ME221-2015-16-(2)"
$ws.Range("B11").Value = "MME235-2017-18-(6)"
$ws.Range("C11").Value = "MME213-2018-19-(7)"
$ws.Range("D11").Value = "MME213-2014-15-(5)"
$ws.Range("B12").Value = "MME235-2018-19-(6)"
$ws.Range("C12").Value = "MME235-2015-16-(1)"
$ws.Range("D12").Value = "MME213-2014-15-(8)"
$ws.Range("B13").Value = "ME221-2017-18-(3)"
$ws.Range("C13").Value = "MME213-2016-17-(3)"
$ws.Range("D13").Value = "EEE267-2018-19-(1)"
$ws.Range("B14").Value = "ME221-2017-18-(8)"
$ws.Range("C14").Value = "ME221-2016-17-(1)"
$ws.Range("D14").Value = "MME213-2015-16-(5)"
$ws.Range("B15").Value = "ME221-2017-18-(4)"
$ws.Range("C15").Value = "ME221-2015-16-(7)"
$ws.Range("D15").Value = "EEE267-2014-15-(3)"
$ws.Range("B16").Value = "MME235-2015-16-(2)"
$ws.Range("C16").Value = "MME213-2018-19-(8)"
$ws.Range("D16").Value = "-"
$ws.Range("B17").Value = "EEE267-2014-15-(7)"
$ws.Range("C17").Value = "MME213-2017-18-(8)"
$ws.Range("D17").Value = "-"
$ws.Range("B18").Value = "MME213-2015-16-(11)"
$ws.Range("C18").Value = "MME213-2018-19-(6)"
$ws.Range("D18").Value = "ME221-2018-19-(2)"
$ws.Range("B19").Value = "EEE267-2014-15-(6)"
$ws.Range("C19").Value = "MME213-2016-17-(1)"
$ws.Range("D19").Value = "EEE267-2016-17-(4)"
$ws.Range("B20").Value = "MME235-2014-15-(5)"
$ws.Range("C20").Value = "MME213-2016-17-(5)"
$ws.Range("D20").Value = "MME235-2017-18-(8)"
$ws.Range("B21").Value = "MME213-2014-15-(3)"
$ws.Range("C21").Value = "EEE267-2015-16-(2)"
$ws.Range("D21").Value = "MME235-2017-18-(4)"
$ws.Range("B22").Value = "MME213-2015-16-(1)"
$ws.Range("C22").Value = "ME221-2018-19-(3)"
$ws.Range("D22").Value = "-"
$ws.Range("B23").Value = "MME235-2016-17-(3)"
$ws.Range("C23").Value = "MME213-2014-15-(7)"
$ws.Range("D23").Value = "ME221-2018-19-(8)"
$ws.Range("B24").Value = "EEE267-2016-17-(7)"
$ws.Range("C24").Value = "EEE267-2016-17-(5)"
$ws.Range("D24").Value = "ME221-2017-18-(6)"
$ws.Range("B25").Value = "ME221-2018-19-(4)"
$ws.Range("C25").Value = "ME221-2015-16-(6)"
$ws.Range("D25").Value = "MME213-2016-17-(12)"
$ws.Range("B26").Value = "EEE267-2015-16-(6)"
$ws.Range("C26").Value = "EEE267-2014-15-(8)"
$ws.Range("D26").Value = "MME235-2014-15-(3)"
$ws.Range("B27").Value = "EEE267-2017-18-(7)"
$ws.Range("C27").Value = "MME213-2016-17-(10)"
$ws.Range("D27").Value = "MME235-2018-19-(2)"
$ws.Range("B28").Value = "ME221-2014-15-(1)"
$ws.Range("C28").Value = "MME213-2015-16-(12)"
$ws.Range("D28").Value = "EEE267-2017-18-(2)"
$ws.Range("B29").Value = "MME213-2016-17-(8)"
$ws.Range("C29").Value = "EEE267-2018-19-(3)"
$ws.Range("D29").Value = "ME221-2014-15-(3)"
$ws.Range("B30").Value = "MME235-2016-17-(8)"
$ws.Range("C30").Value = "EEE267-2018-19-(4)"
$ws.Range("D30").Value = "-"
$ws.Range("B31").Value = "MME213-2018-19-(2)"
$ws.Range("C31").Value = "MME213-2016-17-(7)"
$ws.Range("D31").Value = "EEE267-2015-16-(3)"
$ws.Range("B32").Value = "ME221-2014-15-(6)"
$ws.Range("C32").Value = "EEE267-2014-15-(4)"
$ws.Range("D32").Value = "MME235-2018-19-(4)"
$ws.Range("B33").Value = "EEE267-2017-18-(3)"
$ws.Range("C33").Value = "MME213-2018-19-(1)"
$ws.Range("D33").Value = "MME235-2017-18-(5)"
$ws.Range("B34").Value = "MME235-2014-15-(4)"
$ws.Range("C34").Value = "ME221-2017-18-(5)"
$ws.Range("D34").Value = "MME235-2014-15-(7)"
$ws.Range("B35").Value = "MME235-2015-16-(6)"
$ws.Range("C35").Value = "MME213-2017-18-(5)"
$ws.Range("D35").Value = "MME213-2016-17-(11)"
$ws.Range("B36").Value = "ME221-2017-18-(2)"
$ws.Range("C36").Value = "MME235-2014-15-(8)"
$ws.Range("D36").Value = "MME235-2016-17-(7)"
$ws.Range("B37").Value = "MME213-2017-18-(10)"
$ws.Range("C37").Value = "ME221-2016-17-(5)"
$ws.Range("D37").Value = "EEE267-2014-15-(2)"
$ws.Range("B38").Value = "ME221-2016-17-(2)"
$ws.Range("C38").Value = "MME235-2014-15-(1)"
$ws.Range("D38").Value = "MME213-2015-16-(6)"
$ws.Range("B39").Value = "MME213-2014-15-(10)"
$ws.Range("C39").Value = "ME221-2017-18-(7)"
$ws.Range("D39").Value = "MME235-2017-18-(2)"
$ws.Range("B40").Value = "ME221-2018-19-(7)"
$ws.Range("C40").Value = "MME235-2014-15-(2)"
$ws.Range("D40").Value = "MME213-2015-16-(7)"
$ws.Range("B41").Value = "MME235-2016-17-(5)"
$ws.Range("C41").Value = "EEE267-2016-17-(2)"
$ws.Range("D41").Value = "MME213-2017-18-(3)"
$ws.Range("B42").Value = "MME235-2015-16-(3)"
$ws.Range("C42").Value = "MME235-2015-16-(5)"
$ws.Range("D42").Value = "MME213-2018-19-(3)"
$ws.Range("B43").Value = "MME213-2014-15-(2)"
$ws.Range("C43").Value = "MME213-2017-18-(7)"
$ws.Range("D43").Value = "ME221-2014-15-(5)"
$ws.Range("B44").Value = "EEE267-2017-18-(6)"
$ws.Range("C44").Value = "ME221-2014-15-(2)"
$ws.Range("D44").Value = "ME221-2015-16-(1)"
$ws.Range("B45").Value = "ME221-2014-15-(8)"
$ws.Range("C45").Value = "EEE267-2015-16-(4)"
$ws.Range("D45").Value = "MME213-2017-18-(12)"
$ws.Range("B46").Value = "MME213-2014-15-(1)"
$ws.Range("C46").Value = "MME213-2015-16-(2)"
$ws.Range("D46").Value = "MME213-2014-15-(6)"
$ws.Range("B47").Value = "EEE267-2015-16-(1)"
$ws.Range("C47").Value = "ME221-2016-17-(7)"
$ws.Range("D47").Value = "MME235-2016-17-(4)"
$ws.Range("B48").Value = "ME221-2016-17-(8)"
$ws.Range("C48").Value = "MME235-2018-19-(1)"
$ws.Range("D48").Value = "MME235-2018-19-(3)"
$ws.Range("B49").Value = "MME235-2018-19-(7)"
$ws.Range("C49").Value = "MME213-2017-18-(2)"
$ws.Range("D49").Value = "EEE267-2014-15-(5)"
$ws.Range("B50").Value = "MME235-2018-19-(5)"
$ws.Range("C50").Value = "MME213-2014-15-(11)"
$ws.Range("D50").Value = "MME235-2016-17-(1)"
$ws.Range("B51").Value = "MME213-2018-19-(5)"
$ws.Range("C51").Value = "MME213-2015-16-(10)"
$ws.Range("D51").Value = "MME235-2017-18-(7)"
$ws.Range("B52").Value = "ME221-2018-19-(6)"
$ws.Range("C52").Value = "MME213-2014-15-(9)"
$ws.Range("D52").Value = "MME235-2018-19-(8)"
$ws.Range("B53").Value = "EEE267-2018-19-(7)"
$ws.Range("C53").Value = "MME235-2015-16-(7)"
$ws.Range("D53").Value = "MME235-2015-16-(8)"
$ws.Range("B54").Value = "ME221-2015-16-(4)"
$ws.Range("C54").Value = "MME235-2015-16-(4)"
$ws.Range("D54").Value = "EEE267-2016-17-(6)"
$ws.Range("B55").Value = "MME213-2014-15-(4)"
$ws.Range("C55").Value = "EEE267-2018-19-(5)"
$ws.Range("D55").Value = "EEE267-2018-19-(6)"
$ws.Range("B56").Value = "MME213-2017-18-(9)"
$ws.Range("C56").Value = "MME213-2015-16-(3)"
$ws.Range("D56").Value = "EEE267-2016-17-(8)"
$ws.Range("B57").Value = "ME221-2016-17-(6)"
$ws.Range("C57").Value = "EEE267-2017-18-(5)"
$ws.Range("D57").Value = "MME213-2017-18-(1)"
$ws.Range("B58").Value = "ME221-2017-18-(1)"
$ws.Range("C58").Value = "EEE267-2016-17-(1)"
$ws.Range("D58").Value = "MME213-2017-18-(11)"
$ws.Range("B59").Value = "MME235-2017-18-(1)"
$ws.Range("C59").Value = "ME221-2015-16-(3)"
$ws.Range("D59").Value = "MME213-2016-17-(2)"
$ws.Range("B60").Value = "MME213-2016-17-(4)"
$ws.Range("C60").Value = "MME213-2017-18-(6)"
$ws.Range("D60").Value = "ME221-2018-19-(5)"
$ws.Range("B61").Value = "ME221-2016-17-(4)"
$ws.Range("C61").Value = "MME235-2016-17-(6)"
$ws.Range("D61").Value = "EEE267-2015-16-(7)"
